# "fixed crossover B8 tab in rr_templates"
#
# The rr_crossover_B8 sheet had a stray leading label column (A2:A7,
# "round 1".."round 6", stored as shared strings) that the sibling
# rr_crossover_A8 sheet does not have. This fixes B8 to match A8's
# layout: drop the label column and shift the 8 numeric columns left
# so the data occupies A2:H7 instead of B2:I7.

$wb = $excel.ActiveWorkbook

$wsB8 = $wb.Worksheets.Item("rr_crossover_B8")
$wsA8 = $wb.Worksheets.Item("rr_crossover_A8")

# Remove the leading "round N" label column; Excel shifts B:I left to A:H.
$wsB8.Columns("A:A").Delete()

# rr_crossover_A8 tab: selection moved off the old A2:C4 block to A8.
$wsA8.Range("A8").Select()

# rr_crossover_B8 becomes the active/selected tab, with column A selected.
$wsB8.Activate()
$wsB8.Columns("A:A").Select()
